$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model")

# Change REGID's type from "integer" to "text"
$ws.Range("B10").Value = "text"

# Remove the VISITID row (row 13). This shifts REGIDC/VISITIDC rows up by one,
# turning the former row 14 (REGIDC) into row 13 and former row 15 (VISITIDC)
# into row 14.
$ws.Rows("13").Delete()

# Update the selection to match the final state
$ws.Range("D10").Select()
